# Trade #4 closed at 2026-02-17 19:42:59 - unknown UNKNOWN +0.000%

$wb = $excel.ActiveWorkbook

# --- Sheet: Summary ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B5").Value = 0.45   # Total P&L %
$wsSummary.Range("B6").Value = 4      # Total Trades
$wsSummary.Range("B9").Value = 75     # Win Rate %

# --- Sheet: Strategy Status (MarketMaking row) ---
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("D4").Value = 4       # Trades
$wsStatus.Range("G4").Value = 75      # Win Rate %

# --- Append new trade (#4) to both "All Trades" and "MarketMaking" logs ---
$sheetsToAppend = @("All Trades", "MarketMaking")

foreach ($sheetName in $sheetsToAppend) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Keep the date/time columns as literal text instead of letting Excel
    # auto-convert the date-looking string into a date serial number.
    $ws.Range("B5").NumberFormat = "@"
    $ws.Range("C5").NumberFormat = "@"

    $ws.Range("A5").Value = 4
    $ws.Range("B5").Value = "2026-02-17"
    $ws.Range("C5").Value = "19:42:52"
    $ws.Range("D5").Value = "MarketMaking"
    $ws.Range("E5").Value = "UP"
    $ws.Range("F5").Value = 0.6
    $ws.Range("G5").Value = 0.6
    $ws.Range("H5").Value = "CLOSED"
    $ws.Range("I5").Value = 0
    $ws.Range("J5").Value = 0
    $ws.Range("K5").Value = 100.09
    $ws.Range("L5").Value = 0
    $ws.Range("M5").Value = 0
    $ws.Range("N5").Value = 0.6
    $ws.Range("O5").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P5").Value = "early_exit"
    $ws.Range("Q5").Value = 0.13
}
